# Apply metadata/content updates to the "Metadata" and "Elements" sheets
# per "update to published CDA FHIR logical model with patches #241"

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElem = $wb.Worksheets.Item("Elements")

# --- Metadata sheet updates ---
# Version: 2.0.0-sd-202312-matchbox-patch -> 2.0.0-sd-202406-matchbox-patch
$wsMeta.Range("B3").Value = "2.0.0-sd-202406-matchbox-patch"

# Date: 2024-03-12T18:28:21+01:00 -> 2024-06-19T17:47:42+02:00
$wsMeta.Range("B8").Value = "2024-06-19T17:47:42+02:00"

# Contact: "No display for ContactDetail" -> full HL7 contact string
$wsMeta.Range("B10").Value = "HL7 International - Structured Documents (http://www.hl7.org/Special/committees/structure, structdog@lists.HL7.org)"

# --- Elements sheet updates ---
# Binding Value Set for IVL_TS.operator row (row 5, column Z)
$wsElem.Range("Z5").Value = "http://hl7.org/cda/stds/core/ValueSet/CDASetOperator"

# Column Z width grows to fit the new, longer value (from 49.5 to 51.21484375,
# the closest the COM ColumnWidth setter -- which snaps to pixel granularity -- can reach)
$wsElem.Columns.Item(26).ColumnWidth = 50.3
